$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("loginToForest")
$wsEmail = $wb.Worksheets.Item("email_solution")

# --- loginToForest (sheet6.xml) changes ---

# B2 value changes from "Y" to "N"
$wsLogin.Range("B2").Value = "N"

# Apply wrap-text formatting sheet-wide (matches the default <cols> style=3
# seen in the target) without disturbing already-styled cells (A1 already
# wraps, A2 keeps its Hyperlink+wrap style).
$wsLogin.Cells.WrapText = $true

# New row 3: A3 is a hyperlink (same pattern as A2), B3 = "Y"
$forestUrl = "https://admin.upnorway.com/39042/data/1571867/index/record/1571865/389/summary"
$wsLogin.Hyperlinks.Add($wsLogin.Range("A3"), $forestUrl)
$wsLogin.Range("A3").Style = "Hyperlink"
$wsLogin.Range("A3").WrapText = $true

$wsLogin.Range("B3").Value = "Y"
$wsLogin.Range("B3").WrapText = $true

$wsLogin.Rows.Item(3).RowHeight = 144

# --- email_solution (sheet7.xml): selection moves, no longer the active tab ---
$wsEmail.Activate()
$wsEmail.Range("B5").Select()

# --- loginToForest becomes the active/selected sheet & cell ---
$wsLogin.Activate()
$wsLogin.Range("A4").Select()
